$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string (e.g. "587.35")
# must be forced to remain TEXT (matching the source inline-string cells),
# otherwise Excel auto-converts them to numbers on assignment.
function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
}

$ws.Range('D2').Value = '66.455.59'
$ws.Range('E2').Value = '  -0.69%  '
$ws.Range('D3').Value = '3.467.49'
$ws.Range('E3').Value = '  -1.56%  '
$ws.Range('E4').Value = '  +0.02%  '
Set-TextValue 'D5' '587.35'
$ws.Range('E5').Value = '  +0.23%  '
Set-TextValue 'D6' '176.90'
$ws.Range('E6').Value = '  -0.52%  '
Set-TextValue 'D7' '0.623'
$ws.Range('E7').Value = '  +3.05%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '3.466.28'
$ws.Range('E9').Value = '  -1.53%  '
$ws.Range('E10').Value = '  -1.22%  '
Set-TextValue 'D11' '6.96'
$ws.Range('E11').Value = '  +0.53%  '
Set-TextValue 'D12' '0.417'
$ws.Range('E12').Value = '  -1.92%  '
$ws.Range('D13').Value = '4.068.13'
$ws.Range('E13').Value = '  -1.57%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D14' '30.38'
$ws.Range('E14').Value = '  -0.62%  '
$ws.Range('B15').Value = 'TRON'
$ws.Range('C15').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D15' '0.134'
$ws.Range('E15').Value = '  +1.22%  '
$ws.Range('D16').Value = '66.341.88'
$ws.Range('E16').Value = '  -0.77%  '
$ws.Range('E17').Value = '  -0.43%  '
$ws.Range('D18').Value = '3.450.58'
$ws.Range('E18').Value = '  -1.99%  '
Set-TextValue 'D19' '5.98'
$ws.Range('E19').Value = '  -1.67%  '
Set-TextValue 'D20' '13.83'
$ws.Range('E20').Value = '  -1.73%  '
Set-TextValue 'D21' '373.30'
$ws.Range('E21').Value = '  -2.42%  '
Set-TextValue 'D22' '7.66'
$ws.Range('E22').Value = '  -2.52%  '
Set-TextValue 'D23' '73.28'
$ws.Range('E24').Value = '  -0.03%  '
Set-TextValue 'D25' '0.0000127'
$ws.Range('E25').Value = '  +5.04%  '
$ws.Range('E26').Value = '  -3.16%  '
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '3.615.51'
$ws.Range('E27').Value = '  -1.47%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D28' '9.97'
$ws.Range('E28').Value = '  +0.59%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D29' '0.178'
$ws.Range('E29').Value = '  +2.89%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 'D30' '0.999'
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D31' '5.92'
$ws.Range('E31').Value = '  -0.09%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D32' '2.01'
$ws.Range('E32').Value = '  -0.89%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D33' '23.72'
$ws.Range('E33').Value = '  -3.55%  '
$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D34' '1.00'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D35' '7.06'
$ws.Range('E35').Value = '  -3.19%  '
$ws.Range('B36').Value = 'Fetch.AI'
$ws.Range('C36').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D36' '1.27'
$ws.Range('E36').Value = '  -6.66%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D37' '1.56'
$ws.Range('E37').Value = '  -2.49%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D38' '161.09'
$ws.Range('E38').Value = '  -0.07%  '
$ws.Range('B39').Value = 'Mantle'
$ws.Range('C39').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D39' '0.887'
$ws.Range('E39').Value = '  -1.14%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D40' '28.33'
$ws.Range('E40').Value = '  -5.94%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D41' '1.82'
$ws.Range('E41').Value = '  +0.76%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D42' '2.61'
$ws.Range('E42').Value = '  +1.31%  '
$ws.Range('D43').Value = '2.784.41'
$ws.Range('E43').Value = '  +1.61%  '
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D44' '4.51'
$ws.Range('E44').Value = '  -0.58%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D45' '6.47'
$ws.Range('E45').Value = '  -2.40%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D46' '0.0695'
$ws.Range('E46').Value = '  -1.67%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D47' '25.44'
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D48' '339.95'
$ws.Range('E48').Value = '  +4.97%  '
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D49' '40.04'
$ws.Range('E49').Value = '  -1.82%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D50' '0.0294'
$ws.Range('E50').Value = '  -1.58%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D51' '0.104'
$ws.Range('E51').Value = '  +0.49%  '
